$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.0316002368927
$ws.Range("B1").Value = 1.179243683815002
$ws.Range("C1").Value = 1.48165762424469
$ws.Range("D1").Value = 3.033771753311157
$ws.Range("E1").Value = 4.281867027282715
